$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# dni (D2, Ana Pérez) : "12345678A" -> "12345678A2"
$ws.Range("D2").Value = "12345678A2"

# telefono (F2, Ana Pérez) : "600111222" -> "9"
# Force text format first so the numeric-looking string is NOT
# reinterpreted as a number (keeps it a shared-string / text cell).
$ws.Range("F2").NumberFormat = "@"
$ws.Range("F2").Value = "9"

# telefono (F3, Bruno Díaz) : "600333444" -> "678678"
$ws.Range("F3").NumberFormat = "@"
$ws.Range("F3").Value = "678678"

# direccion (G2, Ana Pérez) : "Calle Sol" -> "Calle Sol78"
$ws.Range("G2").Value = "Calle Sol78"
